$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.277.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.89%  "
$ws.Range("D3").Value = "'1.915.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'327.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").Value = "'0.3946"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.25%  "
$ws.Range("E9").Value = "  +1.77%  "
$ws.Range("D10").Value = "'0.07942"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "
$ws.Range("D11").Value = "'1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'22.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "'1.926.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'7.108"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'5.772"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'0.06955"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'88.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'1.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'0.00001010"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").Value = "'17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "'29.281.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").Value = "'5.361"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "'11.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'2.149.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'2.066"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.82%  "
$ws.Range("D27").Value = "'157.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.33%  "
$ws.Range("D28").Value = "'19.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "
$ws.Range("D29").Value = "'6.113"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.88%  "
$ws.Range("D30").Value = "'1.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "'119.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'0.09396"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").Value = "'0.9267"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("E34").Value = "  +1.03%  "
$ws.Range("D35").Value = "'1.358"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").Value = "'1.204"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").Value = "'0.05842"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("D40").Value = "'7.981"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "'0.5754"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").Value = "'0.1804"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").Value = "'9.987"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").Value = "'2.303"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.86%  "
$ws.Range("D46").Value = "'12.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'0.5430"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("D48").Value = "'0.07076"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("D49").Value = "'1.882"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.05%  "
$ws.Range("D50").Value = "'2.572"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.75%  "
$ws.Range("D51").Value = "'113.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
